# Transition rule summary tables: add "Within 5 miles" and "Within 10 miles"
# of HFC production facility columns to both the Means and Standard
# Deviations sheets, and refresh existing values that changed as a result
# of the script re-run (rows 9 and 10 on both sheets).

$wb = $excel.ActiveWorkbook
$wsMeans = $wb.Worksheets.Item(1)
$wsSD    = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1) Means sheet - new headers first so the new shared strings are
#    created in the same order they appear in the target workbook.
# ---------------------------------------------------------------------
$wsMeans.Cells.Item(1, 6).Value = "Within 5 miles of HFC production facility"
$wsMeans.Cells.Item(1, 7).Value = "Within 10 miles of HFC production facility"

# New column F (Within 5 miles) and G (Within 10 miles) values, rows 2-8
$wsMeans.Cells.Item(2, 6).Value = 50
$wsMeans.Cells.Item(2, 7).Value = 68

$wsMeans.Cells.Item(3, 6).Value = 38
$wsMeans.Cells.Item(3, 7).Value = 23

$wsMeans.Cells.Item(4, 6).Value = 11
$wsMeans.Cells.Item(4, 7).Value = 9.2

$wsMeans.Cells.Item(5, 6).Value = 34
$wsMeans.Cells.Item(5, 7).Value = 25

$wsMeans.Cells.Item(6, 6).Value = 38
$wsMeans.Cells.Item(6, 7).Value = 50

$wsMeans.Cells.Item(7, 6).Value = 11
$wsMeans.Cells.Item(7, 7).Value = 8.6

$wsMeans.Cells.Item(8, 6).Value = 12
$wsMeans.Cells.Item(8, 7).Value = 7.7

# Row 9 (Total Cancer Risk) - values refreshed plus new columns
$wsMeans.Cells.Item(9, 2).Value = 29
$wsMeans.Cells.Item(9, 3).Value = 31
$wsMeans.Cells.Item(9, 4).Value = 41
$wsMeans.Cells.Item(9, 5).Value = 42
$wsMeans.Cells.Item(9, 6).Value = 50
$wsMeans.Cells.Item(9, 7).Value = 59

# Row 10 (Total Respiratory hazard quotient) - values refreshed plus new columns
$wsMeans.Cells.Item(10, 2).Value = 0.37
$wsMeans.Cells.Item(10, 3).Value = 0.36
$wsMeans.Cells.Item(10, 4).Value = 0.4
$wsMeans.Cells.Item(10, 5).Value = 0.4
$wsMeans.Cells.Item(10, 6).Value = 0.4
$wsMeans.Cells.Item(10, 7).Value = 0.4

# ---------------------------------------------------------------------
# 2) Standard Deviations sheet - fill in all numeric data first so that
#    the new "...SD" header strings get appended at the end of the
#    shared string table, matching the target workbook ordering.
# ---------------------------------------------------------------------

# New column F (Within 5 miles SD) and G (Within 10 miles SD) values, rows 2-8
$wsSD.Cells.Item(2, 6).Value = 31
$wsSD.Cells.Item(2, 7).Value = 30

$wsSD.Cells.Item(3, 6).Value = 34
$wsSD.Cells.Item(3, 7).Value = 31

$wsSD.Cells.Item(4, 6).Value = 9.4
$wsSD.Cells.Item(4, 7).Value = 9.4

$wsSD.Cells.Item(5, 6).Value = 26
$wsSD.Cells.Item(5, 7).Value = 23

$wsSD.Cells.Item(6, 6).Value = 13
$wsSD.Cells.Item(6, 7).Value = 25

$wsSD.Cells.Item(7, 6).Value = 9.3
$wsSD.Cells.Item(7, 7).Value = 9.6

$wsSD.Cells.Item(8, 6).Value = 11
$wsSD.Cells.Item(8, 7).Value = 9.8

# Row 9 (Total Cancer Risk SD) - this row previously held "Total Respiratory"
# values and now becomes the "Total Cancer Risk" row; reuses the existing
# shared string created on the Means sheet.
$wsSD.Cells.Item(9, 1).Value = $wsMeans.Cells.Item(9, 1).Value2
$wsSD.Cells.Item(9, 2).Value = 10
$wsSD.Cells.Item(9, 3).Value = 14
$wsSD.Cells.Item(9, 4).Value = 6.3
$wsSD.Cells.Item(9, 5).Value = 4.9
$wsSD.Cells.Item(9, 6).Value = 12
$wsSD.Cells.Item(9, 7).Value = 21

# Row 10 (Total Respiratory hazard quotient SD) - values refreshed plus new columns
$wsSD.Cells.Item(10, 2).Value = 0.14
$wsSD.Cells.Item(10, 3).Value = 0.079
$wsSD.Cells.Item(10, 4).Value = 0.000000000000000025
$wsSD.Cells.Item(10, 5).Value = 0.00000000000000002
$wsSD.Cells.Item(10, 6).Value = 0.015
$wsSD.Cells.Item(10, 7).Value = 0.0097

# New headers added last, so the new shared strings land at the end of
# the shared string table.
$wsSD.Cells.Item(1, 6).Value = "Within 5 mile of HFC production facility SD"
$wsSD.Cells.Item(1, 7).Value = "Within 10 mile of HFC production facility SD"

$wb.Save()
